$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.117.24"
$ws.Range("E2").Value = "  +4.86%  "
$ws.Range("D3").Value = "3.238.78"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.32"
$ws.Range("E5").Value = "  +3.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.80"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "3.235.98"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.68"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +3.36%  "
$ws.Range("D13").Value = "3.801.13"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.77"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "67.046.66"
$ws.Range("E16").Value = "  +4.51%  "
$ws.Range("E17").Value = "  +4.09%  "
$ws.Range("D18").Value = "3.240.91"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.80"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.80"
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.13"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "3.382.59"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.71"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  +4.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.59"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.42"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "174.08"
$ws.Range("E35").Value = "  +10.38%  "
$ws.Range("E36").Value = "  +3.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.75"
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("E38").Value = "  +5.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.847"
$ws.Range("E39").Value = "  +5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  +10.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.59"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("D43").Value = "2.706.45"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.37"
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.27"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.40"
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0669"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.49"
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "332.87"
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("E51").Value = "  +2.62%  "